$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- aeon's ergo buffnerf: adjust ergonomics (C) and weight (D) on existing barrels ---
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 0.54

$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0.63

$ws.Range("C5").Value = -3
$ws.Range("D5").Value = 0.79

$ws.Range("C6").Value = -12
$ws.Range("D6").Value = 0.98

# --- move weight from the bolt carrier group: new row for AUG A3 Std. bolt carrier ---
$ws.Range("A8").Value = "aug_a3_std_bolt_carrier"
$ws.Range("B8").Value = "AUG A3 Std."
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0.2
$ws.Range("M8").Value = 0

# --- view state ---
$excel.ActiveWindow.Zoom = 130
$ws.Range("G11").Select()
